# Edit script: rename the "ZefaniaBibleWithoutCopyright" boolean column into
# a "ZefaniaBibleFreeToEditLicenseType" license-type column.
#
# - Column header H1 is renamed.
# - Existing "False" values become "None".
# - Existing "True" values become either "None" or "PublicDomain" depending on
#   the actual (researched) license of each Bible edition.
# - Column H width is widened to fit the new, longer header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header for column H.
$ws.Range("H1").Value = "ZefaniaBibleFreeToEditLicenseType"

# 2. Rows whose original value was "True" but whose correct license type is
#    "None" (i.e. not actually free to edit / no clear public domain status).
$noneRows = @(7, 36, 38, 62, 63, 74, 75)

# 3. Rows whose original value was "True" and whose correct license type is
#    "PublicDomain".
$publicDomainRows = @(11, 15, 16, 48, 52, 64, 65, 66, 67, 68, 69, 94, 95, 99, 101, 102, 105, 106, 111, 113, 115, 117, 121, 123, 124, 125, 129, 131, 132, 133, 134, 135, 136, 137, 147, 154, 158, 161, 164, 165, 185, 188, 196, 199, 207, 209, 212, 214, 226)

$noneSet = New-Object 'System.Collections.Generic.HashSet[int]'
foreach ($r in $noneRows) { [void]$noneSet.Add($r) }

$pdSet = New-Object 'System.Collections.Generic.HashSet[int]'
foreach ($r in $publicDomainRows) { [void]$pdSet.Add($r) }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()
if ($lastRow -lt 2) { $lastRow = 239 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cur = $cell.Value()

    if ($cur -eq "False") {
        $cell.Value = "None"
    }
    elseif ($cur -eq "True") {
        if ($pdSet.Contains($r)) {
            $cell.Value = "PublicDomain"
        }
        elseif ($noneSet.Contains($r)) {
            $cell.Value = "None"
        }
        else {
            # Fallback (should not happen): default to None.
            $cell.Value = "None"
        }
    }
}

# 4. Widen column H to fit the new header text.
$ws.Columns.Item(8).ColumnWidth = 32.15
